$d = $word.ActiveDocument

# Helper: replace the text found by Find.Execute (searched starting at $SearchRange,
# which is advanced/continued across calls) with $NewText, re-emitting a single
# <w:r> (with optional run-properties XML $RPrXml) via InsertXML so that sibling
# runs in the same paragraph (e.g. a leading empty <w:r/>) are left untouched -
# plain Range.Text / Find-replace rebuilds (merges) the whole paragraph's runs.
function Replace-ExactRun($SearchRange, $OldText, $NewText, $RPrXml) {
    $found = $SearchRange.Find.Execute($OldText)
    if (-not $found) {
        throw "Text not found: $OldText"
    }

    $start = $SearchRange.Start
    $end = $SearchRange.End
    $target = $d.Range($start, $end)

    $escaped = $NewText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:r>' + $RPrXml + '<w:t>' + $escaped + '</w:t></w:r></w:p>' +
           '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xml)

    # Re-seat the search range right after the freshly inserted text so the next
    # Find.Execute on the same $SearchRange continues forward from here.
    $SearchRange.Start = $start
    $SearchRange.End = $start + $NewText.Length
    $SearchRange.Collapse(0)
}

$search = $d.Content

# Title (Heading 1, first occurrence - plain run, no sibling empty run)
Replace-ExactRun $search `
    "Play Eye of Horus Jackpot King for Free - Review" `
    "Play Eye of Horus Jackpot King Free | Review" `
    ""

# "What we like" bullets
Replace-ExactRun $search `
    "Wide betting range" `
    "Straightforward gameplay ideal for beginners" `
    ""

Replace-ExactRun $search `
    "Free play available" `
    "Well-crafted graphics with an ancient Egypt theme" `
    ""

Replace-ExactRun $search `
    "Visually appealing graphics" `
    "Wide range of betting options for any budget" `
    ""

Replace-ExactRun $search `
    "Straightforward gameplay for beginners" `
    "Free gameplay available for testing" `
    ""

# "What we don't like" bullet
Replace-ExactRun $search `
    "High volatility may be off-putting for some players" `
    "High volatility may deter some players" `
    ""

# Title (bold run, second occurrence)
Replace-ExactRun $search `
    "Play Eye of Horus Jackpot King for Free - Review" `
    "Play Eye of Horus Jackpot King Free | Review" `
    "<w:rPr><w:b/></w:rPr>"

# Italic summary paragraph
Replace-ExactRun $search `
    "Eye of Horus Jackpot King review: wide betting range, free play option, and visually appealing graphics. However, below-average RTP and high volatility." `
    "Discover the gameplay, graphics, and betting options of Eye of Horus Jackpot King in our review. Play for free!" `
    "<w:rPr><w:i/></w:rPr>"
